$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$tr.Text = "Power BI Quick Start #2 E04"
$runs = $tr.Runs()
Write-Host "Runs.Count:" $runs.Count
$r1 = $tr.Runs(1)
Write-Host "Run1 text=[$($r1.Text)] start=$($r1.Start)"
